$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = Recorded By
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $first = $parts[0]
            $last = $parts[$parts.Count - 1]
            $parts[0] = $last
            $parts[$parts.Count - 1] = $first
            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
